# Reproduce the "ordenarExcelPorColumna" (buggy) row rotation applied to rows 3-7.
# The B:G data blocks of rows 3-7 are cyclically rotated down by one
# (row7 -> row3, row3 -> row4, row4 -> row5, row5 -> row6, row6 -> row7),
# while column A keeps shifting down by one row too, except row 3's date
# (45460) is left untouched (duplicated into row 4) - this reproduces the
# historical "equal data" sort bug described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") values for rows 3-7, columns A-G.
# NOTE: reading requires the method-call form .Value() in this runtime;
# plain .Value (no parens) does not resolve to the cell's contents.
$before = @{}
foreach ($r in 3..7) {
    $before[$r] = @{
        A = $ws.Cells.Item($r, 1).Value()
        B = $ws.Cells.Item($r, 2).Value()
        C = $ws.Cells.Item($r, 3).Value()
        D = $ws.Cells.Item($r, 4).Value()
        E = $ws.Cells.Item($r, 5).Value()
        F = $ws.Cells.Item($r, 6).Value()
        G = $ws.Cells.Item($r, 7).Value()
    }
}

# Column A (date) target values: row3 stays as-is, rows 4-7 take the
# previous row's original date value.
$ws.Cells.Item(4, 1).Value = $before[3].A
$ws.Cells.Item(5, 1).Value = $before[4].A
$ws.Cells.Item(6, 1).Value = $before[5].A
$ws.Cells.Item(7, 1).Value = $before[6].A
# Row 3's own date (45460) is unchanged.

# Columns B-G: rotate the whole block down by one row, wrapping row 7 -> row 3.
$ws.Cells.Item(3, 2).Value = $before[7].B
$ws.Cells.Item(3, 3).Value = $before[7].C
$ws.Cells.Item(3, 4).Value = $before[7].D
$ws.Cells.Item(3, 5).Value = $before[7].E
$ws.Cells.Item(3, 6).Value = $before[7].F
$ws.Cells.Item(3, 7).Value = $before[7].G

$ws.Cells.Item(4, 2).Value = $before[3].B
$ws.Cells.Item(4, 3).Value = $before[3].C
$ws.Cells.Item(4, 4).Value = $before[3].D
$ws.Cells.Item(4, 5).Value = $before[3].E
$ws.Cells.Item(4, 6).Value = $before[3].F
$ws.Cells.Item(4, 7).Value = $before[3].G

$ws.Cells.Item(5, 2).Value = $before[4].B
$ws.Cells.Item(5, 3).Value = $before[4].C
$ws.Cells.Item(5, 4).Value = $before[4].D
$ws.Cells.Item(5, 5).Value = $before[4].E
$ws.Cells.Item(5, 6).Value = $before[4].F
$ws.Cells.Item(5, 7).Value = $before[4].G

$ws.Cells.Item(6, 2).Value = $before[5].B
$ws.Cells.Item(6, 3).Value = $before[5].C
$ws.Cells.Item(6, 4).Value = $before[5].D
$ws.Cells.Item(6, 5).Value = $before[5].E
$ws.Cells.Item(6, 6).Value = $before[5].F
$ws.Cells.Item(6, 7).Value = $before[5].G

$ws.Cells.Item(7, 2).Value = $before[6].B
$ws.Cells.Item(7, 3).Value = $before[6].C
$ws.Cells.Item(7, 4).Value = $before[6].D
$ws.Cells.Item(7, 5).Value = $before[6].E
$ws.Cells.Item(7, 6).Value = $before[6].F
$ws.Cells.Item(7, 7).Value = $before[6].G
